$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "Base Schema*" text lives in shape named "Subset" (id=43) nested
# inside the group "Group 8" (the 2nd top-level shape on the slide).
$grp = $s.Shapes.Item(2)
$shp = $grp.GroupItems.Item(3)
$shp.TextFrame.TextRange.Text = "Schema Base"
